$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 100; existing rows 100-149 shift down to 101-150.
$ws.Rows("100:100").Insert()

# Populate the newly inserted row 100 with the new observation.
# Columns A-L, Q and T carry over the same categorical metadata as the
# (now shifted-down) row below it.
$ws.Cells.Item(100, 1).Value = 11
$ws.Cells.Item(100, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(100, 3).Value = "Bíobío"
$ws.Cells.Item(100, 4).Value = 45016
$ws.Cells.Item(100, 5).Value = 8
$ws.Cells.Item(100, 6).Value = "Fruta"
$ws.Cells.Item(100, 7).Value = 100101
$ws.Cells.Item(100, 8).Value = "Berries"
$ws.Cells.Item(100, 9).Value = 100101001
$ws.Cells.Item(100, 10).Value = "Arándano (blue)"
$ws.Cells.Item(100, 11).Value = "Sin especificar"
$ws.Cells.Item(100, 12).Value = "Primera"
$ws.Cells.Item(100, 13).Value = 100
$ws.Cells.Item(100, 14).Value = 4000
$ws.Cells.Item(100, 15).Value = 4500
$ws.Cells.Item(100, 16).Value = 4250
$ws.Cells.Item(100, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(100, 18).Value = "Región del Maule"
$ws.Cells.Item(100, 19).Value = 2125
$ws.Cells.Item(100, 20).Value = 2
